$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the confidential disclosure text (shared string) with new "as of" date
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-10 for illustrative purposes only and are subject to change."

# Update the weight/percent-change figures for the holdings rows
$ws.Range("D2").Value = 0.8471312920420058
$ws.Range("E2").Value = -0.01077243236626257

$ws.Range("D3").Value = 0.1528687079579943
$ws.Range("E3").Value = 0.003617945007235779

$ws.Range("E4").Value = -0.008572593970147757
